# Updates "想去人数" (want-to-go count) values in column F across all four
# sheets, matching the gh-pages data refresh captured in the commit.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 67
$ws1.Range("F5").Value = 86
$ws1.Range("F6").Value = 888
$ws1.Range("F7").Value = 481
$ws1.Range("F8").Value = 4797
$ws1.Range("F9").Value = 4797
$ws1.Range("F12").Value = 170
$ws1.Range("F16").Value = 7737
$ws1.Range("F20").Value = 544
$ws1.Range("F21").Value = 1425
$ws1.Range("F23").Value = 6294
$ws1.Range("F24").Value = 2265
$ws1.Range("F27").Value = 1
$ws1.Range("F29").Value = 6222
$ws1.Range("F30").Value = 151
$ws1.Range("F31").Value = 41
$ws1.Range("F32").Value = 120
$ws1.Range("F33").Value = 95
$ws1.Range("F35").Value = 6564
$ws1.Range("F37").Value = 214
$ws1.Range("F38").Value = 101
$ws1.Range("F42").Value = 2477
$ws1.Range("F45").Value = 1118
$ws1.Range("F48").Value = 2167
$ws1.Range("F50").Value = 1097

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 237
$ws2.Range("F4").Value = 46
$ws2.Range("F6").Value = 134
$ws2.Range("F7").Value = 38
$ws2.Range("F8").Value = 11
$ws2.Range("F12").Value = 37

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 1457

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 1457
$ws4.Range("F5").Value = 67
$ws4.Range("F6").Value = 237
$ws4.Range("F7").Value = 86
$ws4.Range("F8").Value = 46
$ws4.Range("F9").Value = 481
$ws4.Range("F10").Value = 4797
$ws4.Range("F11").Value = 4797
$ws4.Range("F14").Value = 170
$ws4.Range("F18").Value = 7737
$ws4.Range("F19").Value = 7737
$ws4.Range("F22").Value = 544
$ws4.Range("F23").Value = 1425
$ws4.Range("F24").Value = 134
$ws4.Range("F25").Value = 6294
$ws4.Range("F26").Value = 2265
$ws4.Range("F29").Value = 6222
$ws4.Range("F30").Value = 151
$ws4.Range("F32").Value = 41
$ws4.Range("F33").Value = 120
$ws4.Range("F34").Value = 95
$ws4.Range("F36").Value = 6564
$ws4.Range("F38").Value = 214
$ws4.Range("F39").Value = 101
$ws4.Range("F42").Value = 37
$ws4.Range("F43").Value = 2477
$ws4.Range("F45").Value = 1118
$ws4.Range("F49").Value = 2167
